$wb = $excel.ActiveWorkbook
$wsVar = $wb.Worksheets.Item("VAR")
$wsValue = $wb.Worksheets.Item("VALUE")

# ---------------------------------------------------------------------------
# Sheet "VAR": add columns B (alternative var name) and C (notations/range)
# ---------------------------------------------------------------------------
$wsVar.Range("B2").Value = "N2"
$wsVar.Range("C2").Value = "0,1,2,3,4,5,6,7,8,9,10"

$wsVar.Range("B3").Value = "N2A"
$wsVar.Range("C3").Value = "0,1,2,3,4,5,6,7,8,9,10"

$wsVar.Range("B4").Value = "F2"
$wsVar.Range("C4").Value = "1,2,3,4"

$wsVar.Range("B5").Value = "F1"
$wsVar.Range("C5").Value = "1,2,3"

$wsVar.Range("B6").Value = "N1B"

$wsVar.Range("B7").Value = "D1"
$wsVar.Range("C7").Value = "1,2,3,4"

$wsVar.Range("B8").Value = "D3"
$wsVar.Range("C8").Value = "1,2,3,4"

$wsVar.Range("B9").Value = "D4"
$wsVar.Range("C9").Value = "1,2,3,4"

$wsVar.Range("B10").Value = "B1"
$wsVar.Range("C10").Value = "1,2,3,4"

$wsVar.Range("B11").Value = "SEX"
$wsVar.Range("C11").Value = "1,2"

$wsVar.Range("B12").Value = "AGE"
$wsVar.Range("C12").Value = "1,2,3,4,5"

$wsVar.Range("B13").Value = "EDU"
$wsVar.Range("C13").Value = "1,2,3,4,5"

$wsVar.Range("B14").Value = "w"

# ---------------------------------------------------------------------------
# Sheet "VALUE": add rows 2-7 describing recode notations
# ---------------------------------------------------------------------------
$wsValue.Range("A2").Value = "d_sup"
$wsValue.Range("B2").Value = "1,2,3"
$wsValue.Range("C2").Value = "1,0,0"
$wsValue.Range("D2").Value = "民主支持編成binary"

$wsValue.Range("A3").Value = "d_sat"
$wsValue.Range("B3").Value = "1,2,3,4"
$wsValue.Range("C3").Value = "4,3,2,1"
$wsValue.Range("D3").Value = "民主滿意反向編碼"

$wsValue.Range("A4").Value = "PartyID"
$wsValue.Range("B4").Value = "1,2,3,4,5,6,7,90,95,96,98,99"
$wsValue.Range("C4").Value = "1,2,1,1,2,0,0,0,0,0,0,0"
$wsValue.Range("D4").Value = "泛藍1泛綠2中立無反應及其他0"

$wsValue.Range("A5").Value = "trust"
$wsValue.Range("B5").Value = "1,2,3,4"
$wsValue.Range("C5").Value = "4,3,2,1"
$wsValue.Range("D5").Value = "政治信任反向編碼"

$wsValue.Range("A6").Value = "engage"
$wsValue.Range("B6").Value = "1,2,3,4"
$wsValue.Range("C6").Value = "4,3,2,1"
$wsValue.Range("D6").Value = "政治參與反向編碼"

$wsValue.Range("A7").Value = "SEX"
$wsValue.Range("B7").Value = "1,2"
$wsValue.Range("C7").Value = "1,0"
$wsValue.Range("D7").Value = "女=0"

# ---------------------------------------------------------------------------
# View state: active sheet switches from VAR to VALUE, with new selections
# ---------------------------------------------------------------------------
$wsVar.Range("D12").Select()
$wsValue.Select()
$wsValue.Range("F8").Select()
